$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.741.59"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "1.726.25"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'240.36"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").Value = "'0.9983"
$ws.Range("D7").Value = "'0.4836"
$ws.Range("E7").Value = "  -0.87%  "
$ws.Range("D8").Value = "'0.2583"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.06185"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "1.728.14"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'15.87"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "'0.06872"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "'0.6041"
$ws.Range("E13").Value = "  +1.06%  "
$ws.Range("D14").Value = "'4.466"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "'77.01"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'0.9981"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "26.563.67"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'0.9974"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'0.000007167"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "'11.36"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "1.946.52"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "'4.417"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "'5.055"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").Value = "'140.11"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "'1.780"
$ws.Range("E27").Value = "  +3.00%  "
$ws.Range("D28").Value = "'106.52"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'1.368"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").Value = "'0.07934"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").Value = "'3.668"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'0.04513"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "'2.595"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "'1.001"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D37").Value = "'0.9358"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("D38").Value = "'2.007"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("D39").Value = "'2.453"
$ws.Range("E39").Value = "  +2.85%  "
$ws.Range("D40").Value = "'0.9975"
$ws.Range("D41").Value = "'0.01498"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'5.614"
$ws.Range("E42").Value = "  +2.57%  "
$ws.Range("D43").Value = "'99.80"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'0.3831"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'6.792"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").Value = "'0.1154"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "'0.05360"
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "'7.942"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").Value = "'51.42"
$ws.Range("E51").Value = "  +0.98%  "
